$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 377; existing rows 377:431 shift down to 378:432,
# carrying their formatting (incl. the date style on column D) with them.
$ws.Rows.Item(377).Insert()

# Populate the newly inserted row 377 with the new weekly price observation.
$ws.Cells.Item(377, 1).Value = 4
$ws.Cells.Item(377, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(377, 3).Value = 'Los Lagos'
$ws.Cells.Item(377, 4).Value = (Get-Date -Year 2022 -Month 8 -Day 3).Date
$ws.Cells.Item(377, 5).Value = 10
$ws.Cells.Item(377, 6).Value = 100114001
$ws.Cells.Item(377, 7).Value = 'Papa'
$ws.Cells.Item(377, 8).Value = 'Patagonia'
$ws.Cells.Item(377, 9).Value = '1a (guarda)'
$ws.Cells.Item(377, 10).Value = 150
$ws.Cells.Item(377, 11).Value = 8000
$ws.Cells.Item(377, 12).Value = 8000
$ws.Cells.Item(377, 13).Value = 8000
$ws.Cells.Item(377, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(377, 15).Value = 'Provincia de Llanquihue'
$ws.Cells.Item(377, 16).Value = 320
$ws.Cells.Item(377, 17).Value = 25
$ws.Cells.Item(377, 18).Value = 'Hortaliza'
